$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" date (column C) for rows 2-7 from 2023-10-13 (45212) to 2023-10-22 (45221)
for ($row = 2; $row -le 7; $row++) {
    $ws.Cells.Item($row, 3).Value = 45221
}
